# Generate Report for Archive
#
# 1. Status text "Ready for handoff" -> "In Translation" wherever it occurs
#    (Status columns on every sheet: Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# 2. Narrow the "Status" column(s):
#      - Overview: columns E (zh-cn) and F (de-de)
#      - zh-cn sheet: column C (Status)
#      - de-de sheet: column C (Status)
#    The stored OOXML width goes from 17.2159881591797 to 13.4101845877511.
#    Excel quantizes ColumnWidth to whole pixels (steps of 1/6 here), so we
#    pick the ColumnWidth value whose resulting stored width lands on the
#    grid point nearest to 13.4101845877511.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newColumnWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            # NOTE: keep the string literal on the LEFT of -eq. Value2 can come
            # back as a .NET boolean for True/False cells, and PowerShell's
            # -eq coerces the right-hand side to the left operand's type -- a
            # boolean left operand would turn any non-empty string into $true
            # and falsely match every text cell.
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value = $newStatus
            }
        }
    }
}

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth   # column F (de-de)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C (Status)
